# Add a "Location" column to the Weekly Data sheet (Tableau-friendly shape)
$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Weekly Data")

# Insert a new column before column B; everything that was in B:G
# (Positive ... Cumulative Deaths) shifts right to C:H, and formulas/
# relative references shift automatically.
$ws.Columns("B").Insert()

# The insert carries column A's formatting into the new column B for every
# row; strip that back off so we can apply clean formatting of our own.
$ws.Range("B1:B43").ClearFormats()

# Drop the placeholder cells the insert created on the still-blank weeks
# (rows 36-43 only have a date in column A, nothing else).
$ws.Range("B36:B43").Clear()

# Header for the new column - match the look of the other header cells.
$ws.Range("B1").Value = "Location"
$ws.Range("A1").Copy()
$ws.Range("B1").PasteSpecial(-4122)

# Fill the new column with the location code for every data row (2-35,
# matching the existing data extent - rows 36-43 are future/blank weeks).
$ws.Range("B2:B35").Value = "pickensC"

# Give the data cells a format, then propagate the resulting style to the
# rest of the column so every row shares a single style entry.
$ws.Range("B2").NumberFormat = "mm-dd-yy"
$ws.Range("B2").Copy()
$ws.Range("B3:B35").PasteSpecial(-4122)
$excel.CutCopyMode = $false

# Column widths: new column B narrower than the data columns, matching the
# original sizing for this location-code helper column.
$ws.Columns("B").ColumnWidth = 8.6640625

# Move the active selection / scroll position to D5 (was F35 / topLeftCell A9)
$ws.Range("D5").Select()
